$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 747.2632
$ws.Range("J17").Value = 747.2632
$ws.Range("L17").Value = 2241.7896
$ws.Range("N17").Value = -2577.7896

$ws.Range("H40").Value = 1642.8572
$ws.Range("I40").Value = 8000
$ws.Range("J40").Value = 1325
$ws.Range("K40").Value = 8000
$ws.Range("L40").Value = 1325
$ws.Range("M40").Value = -7825
$ws.Range("N40").Value = -1675

$ws.Range("H51").Value = 6139.1816
$ws.Range("I51").Value = 5000
$ws.Range("K51").Value = 5000
$ws.Range("M51").Value = -4516

$ws.Range("H113").Value = 4685.7144
$ws.Range("J113").Value = 4125
$ws.Range("L113").Value = 4125
$ws.Range("N113").Value = -10633

$ws.Range("H116").Value = 2433.5
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 2433.5
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 2433.5
$ws.Range("N116").Value = -9317.5
$ws.Range("M116").ClearContents()

$ws.Range("H137").Value = 3336146.2
$ws.Range("I137").Value = 4350042
$ws.Range("J137").Value = 4775.2856
$ws.Range("K137").Value = 13050126
$ws.Range("L137").Value = 14325.8568
$ws.Range("M137").Value = -13047576
$ws.Range("N137").Value = -19425.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7546.6665
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 7546.6665
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 7546.6665
$ws.Range("N32").Value = -8120.6665
$ws.Range("M32").ClearContents()

$ws.Range("H61").Value = 111334390
$ws.Range("I61").Value = 143001150
$ws.Range("J61").Value = 500750
$ws.Range("K61").Value = 143001150
$ws.Range("L61").Value = 500750
$ws.Range("M61").Value = -143000938
$ws.Range("N61").Value = -501174

$ws.Range("H74").Value = 8131531
$ws.Range("I74").Value = 10914620
$ws.Range("J74").Value = 130150
$ws.Range("K74").Value = 10914620
$ws.Range("L74").Value = 130150
$ws.Range("M74").Value = -10913746
$ws.Range("N74").Value = -131898

$ws.Range("H77").Value = 8131531
$ws.Range("I77").Value = 10914620
$ws.Range("J77").Value = 130150
$ws.Range("K77").Value = 54573100
$ws.Range("L77").Value = 650750
$ws.Range("M77").Value = -54568732
$ws.Range("N77").Value = -659486

$ws.Range("H122").Value = 4631807
$ws.Range("I122").Value = 2257.7222
$ws.Range("J122").Value = 18520454
$ws.Range("K122").Value = 6773.1666
$ws.Range("L122").Value = 55561362
$ws.Range("M122").Value = -4323.1666
$ws.Range("N122").Value = -55566262

$ws.Range("H136").Value = 111334390
$ws.Range("I136").Value = 143001150
$ws.Range("J136").Value = 500750
$ws.Range("K136").Value = 429003450
$ws.Range("L136").Value = 1502250
$ws.Range("M136").Value = -429000900
$ws.Range("N136").Value = -1507350

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1102.091
$ws.Range("I99").Value = 1131.1111
$ws.Range("K99").Value = 1131.1111
$ws.Range("M99").Value = 366.8888999999999

$ws.Range("H105").Value = 23811448
$ws.Range("I105").Value = 33335236
$ws.Range("J105").Value = 1983.3334
$ws.Range("K105").Value = 33335236
$ws.Range("L105").Value = 1983.3334
$ws.Range("M105").Value = -33333489
$ws.Range("N105").Value = -5477.3334

$ws.Range("H134").Value = 3314.5
$ws.Range("I134").Value = 3116.681
$ws.Range("J134").Value = 4347.5557
$ws.Range("K134").Value = 9350.043
$ws.Range("L134").Value = 13042.6671
$ws.Range("M134").Value = -6815.043
$ws.Range("N134").Value = -18112.6671

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 200000640
$ws.Range("I22").Value = 200000640
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 200000640
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -200000290
$ws.Range("N22").ClearContents()

$ws.Range("H31").Value = 2901
$ws.Range("I31").Value = 1838.76
$ws.Range("J31").Value = 4463.1177
$ws.Range("K31").Value = 1838.76
$ws.Range("L31").Value = 4463.1177
$ws.Range("M31").Value = -1543.76
$ws.Range("N31").Value = -5053.1177

$ws.Range("H34").Value = 2901
$ws.Range("I34").Value = 1838.76
$ws.Range("J34").Value = 4463.1177
$ws.Range("K34").Value = 1838.76
$ws.Range("L34").Value = 4463.1177
$ws.Range("M34").Value = -1636.76
$ws.Range("N34").Value = -4867.1177

$ws.Range("H58").Value = 43480216
$ws.Range("I58").Value = 55556950
$ws.Range("J58").Value = 3960.2
$ws.Range("K58").Value = 55556950
$ws.Range("L58").Value = 3960.2
$ws.Range("M58").Value = -55556747
$ws.Range("N58").Value = -4366.2

$ws.Range("H132").Value = 62297.65
$ws.Range("I132").Value = 3389.3845
$ws.Range("J132").Value = 253749.5
$ws.Range("K132").Value = 10168.1535
$ws.Range("L132").Value = 761248.5
$ws.Range("M132").Value = -7638.1535
$ws.Range("N132").Value = -766308.5

$ws.Range("H134").Value = 27382.25
$ws.Range("I134").Value = 1832.975
$ws.Range("J134").Value = 282875
$ws.Range("K134").Value = 5498.924999999999
$ws.Range("L134").Value = 848625
$ws.Range("M134").Value = -2963.924999999999
$ws.Range("N134").Value = -853695

$ws.Range("H136").Value = 43480216
$ws.Range("I136").Value = 55556950
$ws.Range("J136").Value = 3960.2
$ws.Range("K136").Value = 166670850
$ws.Range("L136").Value = 11880.6
$ws.Range("M136").Value = -166668300
$ws.Range("N136").Value = -16980.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H19").Value = 817.6
$ws.Range("J19").Value = 817.6
$ws.Range("L19").Value = 2452.8
$ws.Range("N19").Value = -2800.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1448.6666
$ws.Range("I97").Value = 1614.8823
$ws.Range("J97").Value = 742.25
$ws.Range("K97").Value = 1614.8823
$ws.Range("L97").Value = 742.25
$ws.Range("M97").Value = -1118.8823
$ws.Range("N97").Value = -1734.25

$ws.Range("H126").Value = 2283.3333
$ws.Range("I126").Value = 2225.8462
$ws.Range("J126").Value = 2376.75
$ws.Range("K126").Value = 6677.5386
$ws.Range("L126").Value = 7130.25
$ws.Range("M126").Value = -4207.5386
$ws.Range("N126").Value = -12070.25

$ws.Range("H141").Value = 29347.375
$ws.Range("J141").Value = 29347.375
$ws.Range("L141").Value = 29347.375
$ws.Range("N141").Value = -39707.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2578.513
$ws.Range("I16").Value = 1224.2693
$ws.Range("J16").Value = 5287
$ws.Range("K16").Value = 1224.2693
$ws.Range("L16").Value = 5287
$ws.Range("M16").Value = -1054.2693
$ws.Range("N16").Value = -5627

$ws.Range("H82").Value = 5520.3
$ws.Range("I82").Value = 1300
$ws.Range("J82").Value = 5989.222
$ws.Range("K82").Value = 1300
$ws.Range("L82").Value = 5989.222
$ws.Range("M82").Value = -939
$ws.Range("N82").Value = -6711.222

$ws.Range("H85").Value = 5520.3
$ws.Range("I85").Value = 1300
$ws.Range("J85").Value = 5989.222
$ws.Range("K85").Value = 1300
$ws.Range("L85").Value = 5989.222
$ws.Range("M85").Value = -52
$ws.Range("N85").Value = -8485.222

$ws.Range("H132").Value = 218156.86
$ws.Range("I132").Value = 500000
$ws.Range("K132").Value = 1500000
$ws.Range("M132").Value = -1497470

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 3250
$ws.Range("I96").Value = 3250
$ws.Range("K96").Value = 3250
$ws.Range("M96").Value = -1877

$ws.Range("H132").Value = 134845
$ws.Range("I132").Value = 78282.766
$ws.Range("J132").Value = 502499.5
$ws.Range("K132").Value = 234848.298
$ws.Range("L132").Value = 1507498.5
$ws.Range("M132").Value = -232318.298
$ws.Range("N132").Value = -1512558.5

$ws.Range("H136").Value = 51701.875
$ws.Range("I136").Value = 40883
$ws.Range("J136").Value = 69733.336
$ws.Range("K136").Value = 122649
$ws.Range("L136").Value = 209200.008
$ws.Range("M136").Value = -120099
$ws.Range("N136").Value = -214300.008
